$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (swap column pairs)
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Update the one-hot indicator values for each data row
$ws.Range("B2").Value = 1
$ws.Range("F2").Value = 0

$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = 0
$ws.Range("E5").Value = 1

$ws.Range("A6").Value = 0
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = 1
$ws.Range("D7").Value = 0
